# The SSSOM schema workbook used to have a dedicated "EntityReference" sheet.
# That type has been folded into "uriorcurie", so the sheet is no longer
# needed and is removed. The remaining four schema sheets (Mapping,
# MappingRegistry, MappingSet, MappingSetReference) are each duplicated,
# with the duplicates appended at the end named with a trailing "1"
# (Mapping1, MappingRegistry1, MappingSet1, MappingSetReference1).

$wb = $excel.ActiveWorkbook

# Remove the obsolete EntityReference sheet.
$entityRef = $wb.Worksheets.Item("EntityReference")
[void]$entityRef.Delete()

# Duplicate each remaining sheet, appending the copy at the end of the
# workbook and naming it "<OriginalName>1".
$sheetNames = @("Mapping", "MappingRegistry", "MappingSet", "MappingSetReference")
foreach ($sheetName in $sheetNames) {
    $source = $wb.Worksheets.Item($sheetName)
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $source.Copy($null, $lastSheet)

    $copy = $wb.Worksheets.Item($wb.Worksheets.Count)
    $copy.Name = $sheetName + "1"
}

# Leave the final (newly added) sheet as the active one.
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
